# "fully fucntional CO2 objective" -- shrink the df_con_electric table from
# 2 of each asset (net1/net2, pv1/pv2, bat1/bat2, CHP1/CHP2, demand1/demand2)
# down to a single instance of each (net1, pv1, bat1, CHP1, demand1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "...2" columns (F:I) and the "...2" rows (5:7) that are no
# longer part of the (now smaller) dataframe.
$ws.Range("F1:I7").EntireColumn.Delete() | Out-Null
$ws.Range("A5:E7").EntireRow.Delete() | Out-Null

# Relabel the remaining header row / index column to match the new,
# de-duplicated asset list.
$ws.Range("B1").Value = "net1"
$ws.Range("C1").Value = "pv1"
$ws.Range("D1").Value = "bat1"
$ws.Range("E1").Value = "CHP1"

$ws.Range("A2").Value = "demand1"
$ws.Range("A3").Value = "net1"
$ws.Range("A4").Value = "bat1"
